# Add a new progress-log row (row 8) to the "Tiến trình" sheet, matching the
# formatting conventions already used by the table (row 6/7), and move the
# "latest entry" highlight from row 7 to the new row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# 1) Clone the formatting of the current last data row (row 7) onto the new
#    row 8 so fonts/borders/number-formats/wrap stay consistent with the rest
#    of the table.
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

# 2) Row 7 is no longer the newest entry, so its highlight cell reverts back
#    to the plain "already covered" look used by row 6.
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)

# 3) Populate the new row's content.
$ws.Cells.Item(8,1).Value = 46003
$ws.Cells.Item(8,2).Value = "Tầng DAO đã hoàn thiện cơ bản, đã điều chỉnh các model theo cơ chế mới"
$ws.Cells.Item(8,3).Value = "hoàn thiện tầng Service, test login "
$ws.Cells.Item(8,4).Value = "điều chỉnh theo các phương thức mã hoá BCrypt để lưu trữ mật khẩu"

# 4) Match the wrap-driven row height used elsewhere in the sheet.
$ws.Rows.Item(8).RowHeight = 57.6

# 5) Sheet view bookkeeping: zoom level and active selection moved to D6.
$ws.Application.ActiveWindow.Zoom = 98
$ws.Range("D6").Select()
